$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 150; this shifts the existing rows 150:176
# down to 151:177 (and extends the used range to R177).
$ws.Rows(150).Insert()

# Populate the newly inserted row 150 with the new weekly data point
# (a new "Albahaca" price observation for Feria Lagunitas de Puerto Montt).
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C150").Value = "Los Lagos"
$ws.Range("D150").Value = 44995
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = 100112052
$ws.Range("G150").Value = "Albahaca"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 80
$ws.Range("K150").Value = 6500
$ws.Range("L150").Value = 6500
$ws.Range("M150").Value = 6500
$ws.Range("N150").Value = "$/docena de matas"
$ws.Range("O150").Value = "Región Metropolitana"
$ws.Range("P150").Value = 1083
$ws.Range("Q150").Value = 6
$ws.Range("R150").Value = "Hortaliza"
